$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.382.61'
$ws.Range('E2').Value = '  -3.85%  '
$ws.Range('D3').Value = '1.766.50'
$ws.Range('E3').Value = '  -3.13%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').Value = '''305.67'
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('D7').Value = '''0.4302'
$ws.Range('E7').Value = '  +1.41%  '
$ws.Range('D8').Value = '''0.3630'
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('D9').Value = '''0.07083'
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('D10').Value = '''0.8455'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('D11').Value = '''20.28'
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('D12').Value = '1.816.35'
$ws.Range('E12').Value = '  -4.68%  '
$ws.Range('D13').Value = '''5.247'
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('D14').Value = '''6.433'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').Value = '''0.06818'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '''79.16'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('E18').Value = '  -2.70%  '
$ws.Range('D19').Value = '''1.001'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = '''15.03'
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').Value = '26.377.27'
$ws.Range('E21').Value = '  -5.12%  '
$ws.Range('D22').Value = '''5.029'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('D23').Value = '''11.20'
$ws.Range('E23').Value = '  +2.30%  '
$ws.Range('D24').Value = '1.972.75'
$ws.Range('E24').Value = '  -5.47%  '
$ws.Range('D25').Value = '''152.69'
$ws.Range('E25').Value = '  -1.53%  '
$ws.Range('D26').Value = '''1.858'
$ws.Range('E26').Value = '  -6.46%  '
$ws.Range('D27').Value = '''18.16'
$ws.Range('E27').Value = '  -3.08%  '
$ws.Range('D28').Value = '''5.078'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').Value = '''114.10'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '''1.700'
$ws.Range('E30').Value = '  -5.30%  '
$ws.Range('D31').Value = '''0.08924'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').Value = '''0.7297'
$ws.Range('E32').Value = '  -2.63%  '
$ws.Range('D33').Value = '''4.334'
$ws.Range('E33').Value = '  -4.39%  '
$ws.Range('D34').Value = '''1.111'
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('D35').Value = '''2.761'
$ws.Range('E35').Value = '  -7.16%  '
$ws.Range('D36').Value = '''1.001'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '''1.072'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('D38').Value = '''0.05130'
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').Value = '''0.01890'
$ws.Range('E39').Value = '  -1.35%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = '''0.1610'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.4916'
$ws.Range('E41').Value = '  -2.88%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''6.239'
$ws.Range('E42').Value = '  -3.24%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '''2.503'
$ws.Range('E43').Value = '  -9.99%  '
$ws.Range('E44').Value = '  -3.29%  '
$ws.Range('D45').Value = '''105.01'
$ws.Range('E45').Value = '  -0.62%  '
$ws.Range('D46').Value = '''1.001'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '''10.07'
$ws.Range('E47').Value = '  -3.13%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.06186'
$ws.Range('E48').Value = '  -4.01%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = '''0.4479'
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('E51').Value = '  +1.26%  '
